# Regen save_data to use K (constant = 1) instead of Strike# in column G
# for the data rows (row 2 through row 6 of Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
